# 31 mars 2024 v1.0
# Feuille.xlsx template update:
#   - "Metiers" column header renamed to "Statut"
#   - old "Statut" column renamed/replaced by a new "Montant Révisé" column
#   - the H1:M1 header cells are re-formatted to reuse the same (border +
#     centered alignment, no explicit fill) style already used by A1:G1
#   - columns L (12) and M (13) are resized
#   - the sheet view scrolls right and the active cell/selection moves to N6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text changes -----------------------------------------------
# Column L used to read "Metiers" -> now "Statut"
# Column M used to read "Statut" -> now "Montant Révisé"
$ws.Range("L1").Value = "Statut"
$ws.Range("M1").Value = "Montant Révisé"

# --- Re-unify the cell style on H1:M1 with the one used on A1:G1 -------
# (border all around + centered horizontal/vertical alignment, default fill)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("H1:M1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column widths -------------------------------------------------------
# (values picked so the exported <col width=.../> lands on 27.85546875 /
# 16 given this host's char-width -> pixel rounding)
$ws.Columns.Item(12).ColumnWidth = 27
$ws.Columns.Item(13).ColumnWidth = 15.16666666666667

# --- View state: scrolled right, N6 selected -----------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("N6").Select() | Out-Null

Write-Host "Feuille.xlsx header update applied"
